$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New data row: Email + generated password, mirroring the existing
# Email_ID / Password rows above it.
$ws.Range("A3").Value = "syed_hussain@hcl.com"
$ws.Range("B3").Value = "hhhjjj"

# Excel auto-hyperlinks the email in A3 (like A2) but not the
# plain-text password in B3 (like B2 does get one - here it doesn't,
# matching the authored diff).
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:syed_hussain@hcl.com")

# Re-apply the Hyperlink cell style so A3 reuses the same style record
# as A2/B2 instead of a freshly synthesized one.
$ws.Range("A3").Style = "Hyperlink"

# Reflect the cursor move that was captured in the saved view state.
$ws.Range("C10").Select()
